$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 211
$ws1.Range("F5").Value = 1180
$ws1.Range("F7").Value = 568
$ws1.Range("F8").Value = 114
$ws1.Range("F9").Value = 569
$ws1.Range("F10").Value = 603
$ws1.Range("F11").Value = 86
$ws1.Range("F12").Value = 41
$ws1.Range("F13").Value = 148

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 17

# Sheet "本地生活" (local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6275
$ws3.Range("F3").Value = 773
$ws3.Range("F4").Value = 1899

# Sheet "全部类型" (all types) - aggregated view
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6275
$ws4.Range("F3").Value = 773
$ws4.Range("F4").Value = 1899
$ws4.Range("F11").Value = 211
$ws4.Range("F13").Value = 17
$ws4.Range("F15").Value = 1180
$ws4.Range("F19").Value = 568
$ws4.Range("F21").Value = 114
$ws4.Range("F22").Value = 569
$ws4.Range("F24").Value = 603
$ws4.Range("F25").Value = 86
$ws4.Range("F28").Value = 41
$ws4.Range("F29").Value = 148
